# data cleanup continued in player_per_game_df
# Remove "Joel Embiid" and "LeBron James" rows from the player/award count
# table, shifting the remaining rows up (mirrors a pandas dropna/filter
# operation re-exported to the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$playersToRemove = @("Joel Embiid", "LeBron James")

# Collect the row numbers for the players to remove first, then delete
# them from the bottom up so earlier row numbers remain valid targets.
$rowsToDelete = @()
foreach ($player in $playersToRemove) {
    $cell = $ws.Cells.Find($player)
    if ($cell -ne $null) {
        $rowsToDelete += $cell.Row
    }
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
